$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated LR-pair table: 3 sender clusters (ECs, FAPs, MuSCs) x 3 target
# clusters (ECs, FAPs, MuSCs) with refreshed TPM-derived statistics (rows 2-10).
$rows = @{
  2 = @("ECs", "Cntn2", "Cntnap2", "ECs", 2, 0.6666666666666666, 0.05028033333333334, 0.150841, 0.1252715694221136, 0.1252715694221136, 1, 0.3333333333333333, 0.00217, 0.00651, 0.01227840867296996, 0.01227840867296996, 0.0001091083233333333, 0.0009819749100000001, 0.001538135524469038, 0.001538135524469038)
  3 = @("ECs", "Cntn2", "Cntnap2", "FAPs", 2, 0.6666666666666666, 0.05028033333333334, 0.150841, 0.1252715694221136, 0.1252715694221136, 2, 0.6666666666666666, 0.08864699999999999, 0.265941, 0.5015871399229346, 0.5015871399229346, 0.004457200709, 0.040114806381, 0.0628346082200953, 0.06283460822009529)
  4 = @("ECs", "Cntn2", "Cntnap2", "MuSCs", 2, 0.6666666666666666, 0.05028033333333334, 0.150841, 0.1252715694221136, 0.1252715694221136, 3, 1, 0.08591599999999999, 0.257748, 0.4861344514040954, 0.4861344514040954, 0.004319885118666667, 0.038878966068, 0.06089882567754925, 0.06089882567754924)
  5 = @("FAPs", "Cntn2", "Cntnap2", "ECs", 2, 0.6666666666666666, 0.300673, 0.9020189999999999, 0.7491155307811899, 0.7491155307811898, 1, 0.3333333333333333, 0.00217, 0.00651, 0.01227840867296996, 0.01227840867296996, 0.00065246041, 0.00587214369, 0.009197946630200259, 0.009197946630200255)
  6 = @("FAPs", "Cntn2", "Cntnap2", "FAPs", 2, 0.6666666666666666, 0.300673, 0.9020189999999999, 0.7491155307811899, 0.7491155307811898, 2, 0.6666666666666666, 0.08864699999999999, 0.265941, 0.5015871399229346, 0.5015871399229346, 0.026653759431, 0.2398838348789999, 0.3757467165563881, 0.375746716556388)
  7 = @("FAPs", "Cntn2", "Cntnap2", "MuSCs", 2, 0.6666666666666666, 0.300673, 0.9020189999999999, 0.7491155307811899, 0.7491155307811898, 3, 1, 0.08591599999999999, 0.257748, 0.4861344514040954, 0.4861344514040954, 0.025832621468, 0.2324935932119999, 0.3641708675946015, 0.3641708675946014)
  8 = @("MuSCs", "Cntn2", "Cntnap2", "ECs", 1, 0.3333333333333333, 0.05041733333333333, 0.151252, 0.1256128997966967, 0.1256128997966967, 1, 0.3333333333333333, 0.00217, 0.00651, 0.01227840867296996, 0.01227840867296996, 0.0001094056133333333, 0.00098465052, 0.001542326518300667, 0.001542326518300667)
  9 = @("MuSCs", "Cntn2", "Cntnap2", "FAPs", 1, 0.3333333333333333, 0.05041733333333333, 0.151252, 0.1256128997966967, 0.1256128997966967, 2, 0.6666666666666666, 0.08864699999999999, 0.265941, 0.5015871399229346, 0.5015871399229346, 0.004469345348, 0.040224108132, 0.06300581514645126, 0.06300581514645125)
  10 = @("MuSCs", "Cntn2", "Cntnap2", "MuSCs", 1, 0.3333333333333333, 0.05041733333333333, 0.151252, 0.1256128997966967, 0.1256128997966967, 3, 1, 0.08591599999999999, 0.257748, 0.4861344514040954, 0.4861344514040954, 0.004331655610666667, 0.03898490049599999, 0.06106475813194475, 0.06106475813194474)
}

foreach ($r in $rows.Keys) {
  $vals = $rows[$r]
  for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
  }
}
